$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths
$ws.Columns.Item(2).ColumnWidth = 7.16   # B -> raw width 8
$ws.Columns.Item(3).ColumnWidth = 7.16   # C -> raw width 8
$ws.Columns.Item(7).ColumnWidth = 7.16   # G -> raw width 8
$ws.Columns.Item(10).ColumnWidth = 7.16   # J -> raw width 8
$ws.Columns.Item(11).ColumnWidth = 7.16   # K -> raw width 8
$ws.Columns.Item(12).ColumnWidth = 7.16   # L -> raw width 8
$ws.Columns.Item(15).ColumnWidth = 7.16   # O -> raw width 8
$ws.Columns.Item(16).ColumnWidth = 7.16   # P -> raw width 8
$ws.Columns.Item(17).ColumnWidth = 7.16   # Q -> raw width 8
$ws.Columns.Item(20).ColumnWidth = 8.16   # T -> raw width 9
$ws.Columns.Item(22).ColumnWidth = 7.16   # V -> raw width 8
$ws.Columns.Item(27).ColumnWidth = 7.16   # AA -> raw width 8
$ws.Columns.Item(28).ColumnWidth = 7.16   # AB -> raw width 8
$ws.Columns.Item(30).ColumnWidth = 7.16   # AD -> raw width 8
$ws.Columns.Item(34).ColumnWidth = 7.16   # AH -> raw width 8

# Update data values for rows 2-5 (new sensor readings)
$data = New-Object 'object[,]' 4,34
$data[0,0] = 45095.50694444445
$data[0,1] = 7.205
$data[0,2] = 5.605
$data[0,3] = 1.192
$data[0,4] = 15.61
$data[0,5] = 12.374
$data[0,6] = 4.794
$data[0,7] = 14.798
$data[0,8] = 8.957000000000001
$data[0,9] = 4.429
$data[0,10] = 5.631
$data[0,11] = 6.249
$data[0,12] = 7.306
$data[0,13] = 2.788
$data[0,14] = 6.015
$data[0,15] = 7.738
$data[0,16] = 5.138
$data[0,17] = 0.492
$data[0,18] = 0.931
$data[0,19] = 84.14
$data[0,20] = 16.424
$data[0,21] = 5.552
$data[0,22] = 10.174
$data[0,23] = 6.283
$data[0,24] = 0.894
$data[0,25] = 9.741
$data[0,26] = 4.361
$data[0,27] = 4.885
$data[0,28] = 6.06
$data[0,29] = 8.134
$data[0,30] = 1.522
$data[0,31] = 13.18
$data[0,32] = 3.625
$data[0,33] = 6.387
$data[1,0] = 45095.51388888889
$data[1,1] = 20.379
$data[1,2] = 15.362
$data[1,3] = 1.081
$data[1,4] = 44.569
$data[1,5] = 36.329
$data[1,6] = 15.496
$data[1,7] = 57.659
$data[1,8] = 24.785
$data[1,9] = 11.565
$data[1,10] = 16.316
$data[1,11] = 17.829
$data[1,12] = 19.25
$data[1,13] = 5.716
$data[1,14] = 16.166
$data[1,15] = 22.637
$data[1,16] = 13.594
$data[1,17] = 0.32
$data[1,18] = 0.958
$data[1,19] = 238.743
$data[1,20] = 45.131
$data[1,21] = 14.922
$data[1,22] = 30.052
$data[1,23] = 16.29
$data[1,24] = 2.149
$data[1,25] = 30.233
$data[1,26] = 12.881
$data[1,27] = 11.965
$data[1,28] = 14.16
$data[1,29] = 19.444
$data[1,30] = 0.707
$data[1,31] = 52.871
$data[1,32] = 8.786
$data[1,33] = 18.354
$data[2,0] = 45095.52083333334
$data[2,1] = 7.01
$data[2,2] = 5.282
$data[2,3] = 0.487
$data[2,4] = 15.484
$data[2,5] = 12.38
$data[2,6] = 5.124
$data[2,7] = 25.999
$data[2,8] = 8.555
$data[2,9] = 4.232
$data[2,10] = 5.471
$data[2,11] = 6.17
$data[2,12] = 6.793
$data[2,13] = 2.193
$data[2,14] = 5.639
$data[2,15] = 7.773
$data[2,16] = 4.848
$data[2,17] = 0.135
$data[2,18] = 0.442
$data[2,19] = 78.544
$data[2,20] = 15.979
$data[2,21] = 5.205
$data[2,22] = 10.471
$data[2,23] = 5.798
$data[2,24] = 0.753
$data[2,25] = 13.07
$data[2,26] = 4.392
$data[2,27] = 4.344
$data[2,28] = 5.133
$data[2,29] = 6.914
$data[2,30] = 0.461
$data[2,31] = 24.325
$data[2,32] = 3.164
$data[2,33] = 6.302
$data[3,0] = 45095.52777777778
$data[3,1] = 12.82
$data[3,2] = 9.630000000000001
$data[3,3] = 0.61
$data[3,4] = 28.1
$data[3,5] = 22.88
$data[3,6] = 9.779999999999999
$data[3,7] = 38.19
$data[3,8] = 15.57
$data[3,9] = 7.28
$data[3,10] = 10.21
$data[3,11] = 11.24
$data[3,12] = 12.09
$data[3,13] = 3.55
$data[3,14] = 10.15
$data[3,15] = 14.23
$data[3,16] = 8.56
$data[3,17] = 0.15
$data[3,18] = 0.55
$data[3,19] = 147.17
$data[3,20] = 28.33
$data[3,21] = 9.369999999999999
$data[3,22] = 18.9
$data[3,23] = 10.22
$data[3,24] = 1.34
$data[3,25] = 19.26
$data[3,26] = 8.119999999999999
$data[3,27] = 7.51
$data[3,28] = 8.84
$data[3,29] = 12.12
$data[3,30] = 0.34
$data[3,31] = 34.76
$data[3,32] = 5.5
$data[3,33] = 11.56
$ws.Range("A2:AH5").Value2 = $data

# Remove the now-unused last data row (row 6); this also fixes the sheet dimension
$ws.Rows.Item(6).Delete()
